# disable mousewheel on comboboxes
#
# - Add a new "Todo"/"Feature" task to the Active sheet: id 23,
#   "when combox value is changed, return focus to textbox", created 3/3/2018.
# - The two related tasks that are now finished move from Active to the
#   Inactive (Done) sheet, with a Done date of 3/3/2018:
#     id 9  "deactivate mouse scroll when hovering over dropdowns..."
#     id 10 "icon for program"
# - Bump the Config sheet's Max Id to 23.

$wb = $excel.ActiveWorkbook
$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")
$config = $wb.Worksheets.Item("Config")

# Helper pattern used throughout: text that looks like a date (e.g.
# "3/3/2018") gets auto-converted to a date serial by a plain .Value
# assignment, so force text storage with NumberFormat "@" first, then put
# the cell style back to Normal (General) once the text value is in place.

# --- Active sheet: insert the new task as row 2 (just under the header) ---
$active.Rows.Item(2).Insert()
$active.Range("A2:F2").Style = "Normal"

$active.Range("A2").Value = 23
$active.Range("B2").Value = "when combox value is changed, return focus to textbox"
$active.Range("C2").Value = "Todo"
$active.Range("D2").Value = "Feature"
$active.Range("E2").NumberFormat = "@"
$active.Range("E2").Value = "3/3/2018"
$active.Range("E2").Style = "Normal"
$active.Range("F2").ClearContents()

# --- Active sheet: remove the two tasks that are now done ---
# After the insert above, former row 4 (id 9) and row 5 (id 10) are now at
# row 5 and row 6. Capture their title text before deleting, then move
# them over to the Inactive sheet.
$deactivateTitle = $active.Range("B5").Value()
$iconTitle = $active.Range("B6").Value()

$active.Rows.Item(6).Delete()
$active.Rows.Item(5).Delete()

# --- Inactive sheet: insert the two finished tasks at the top of the list ---
$inactive.Rows.Item(2).Insert()
$inactive.Rows.Item(2).Insert()
$inactive.Range("A2:F3").Style = "Normal"

$inactive.Range("A2").Value = 10
$inactive.Range("B2").Value = $iconTitle
$inactive.Range("C2").Value = "Done"
$inactive.Range("D2").Value = "Feature"
$inactive.Range("E2").NumberFormat = "@"
$inactive.Range("E2").Value = "12/5/2017"
$inactive.Range("E2").Style = "Normal"
$inactive.Range("F2").NumberFormat = "@"
$inactive.Range("F2").Value = "3/3/2018"
$inactive.Range("F2").Style = "Normal"

$inactive.Range("A3").Value = 9
$inactive.Range("B3").Value = $deactivateTitle
$inactive.Range("C3").Value = "Done"
$inactive.Range("D3").Value = "Bug"
$inactive.Range("E3").NumberFormat = "@"
$inactive.Range("E3").Value = "12/1/2017"
$inactive.Range("E3").Style = "Normal"
$inactive.Range("F3").NumberFormat = "@"
$inactive.Range("F3").Value = "3/3/2018"
$inactive.Range("F3").Style = "Normal"

# --- Config sheet: bump Max Id to match the newly created task ---
$config.Range("F2").Value = 23
